$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(1.032923005835927, 1.033758375452343, 1.034299558342161, 1.034527223911973, 1.034565458845671, 1.03430259982431, 1.033205188816147, 1.031276406169549, 1.029994006856443, 1.029439553845574, 1.029233732297491, 1.029277875945509, 1.029422537972234, 1.029511685897351, 1.03003082179508, 1.030356686561599, 1.030546838359964, 1.030611688795648, 1.030321716028858, 1.029379935092354, 1.028788535559078, 1.029101977375543, 1.03033751745911, 1.031774440720077)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @(1.040505655009929, 1.041166811567473, 1.041595200698994, 1.041775431780621, 1.04180570131338, 1.041597608422205, 1.040728975514282, 1.039202827242327, 1.03818852128558, 1.037750079261569, 1.037587338221161, 1.037622241455547, 1.037736624662741, 1.037807115336602, 1.038217635386422, 1.038475348281438, 1.038625740979454, 1.038677033351861, 1.038447690550003, 1.037702938446739, 1.037235354394219, 1.037483165330666, 1.03846018766261, 1.039596829801183)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(1.032304413604033, 1.03301145641942, 1.033469892948289, 1.033662841115578, 1.033695250895382, 1.033472470263262, 1.032543167859547, 1.030912840235898, 1.029830918467643, 1.029363634411482, 1.029190245725779, 1.029227429948159, 1.029349298337047, 1.029424409568688, 1.029861955911294, 1.030136738471521, 1.030297129686203, 1.030351838428442, 1.030107244974677, 1.029313406111105, 1.028815338657466, 1.029079273440761, 1.030120571460674, 1.031333452065745)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(1.049205897372537, 1.050099456979084, 1.050678777134007, 1.050922591310489, 1.050963544442911, 1.050682033940361, 1.049507645457044, 1.047446941243915, 1.046079113820194, 1.045488270718632, 1.045269022847341, 1.045316042383575, 1.045470143171726, 1.045565118640703, 1.046118356534634, 1.046465773662758, 1.046668554541383, 1.046737720987442, 1.046428484769091, 1.045424758309048, 1.044794935365773, 1.045128696369043, 1.046445333587837, 1.04797863780388)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colI = @(1.037130258665913, 1.037314973998913, 1.037433622975956, 1.037483293026659, 1.037491620513276, 1.037434287495544, 1.037192864666845, 1.036760782320044, 1.03646828817839, 1.036340590818621, 1.036293001956368, 1.036303216994122, 1.036336660294455, 1.036357245114822, 1.036476741063878, 1.036551418302783, 1.036594875353916, 1.036609675946694, 1.03654341657531, 1.036326816387262, 1.036189726846731, 1.036262486101428, 1.036547032523433, 1.036873271462434)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $colI[$i]
}

$colJ = @(1.038050560673331, 1.038529007560387, 1.038838576393823, 1.038968713498941, 1.038990563724829, 1.038840315315892, 1.038212257073823, 1.037105453497297, 1.036367602216783, 1.036048124442817, 1.035929459948441, 1.035954913694104, 1.036038315514345, 1.03608970266661, 1.03638880536442, 1.036576429949538, 1.036685869721969, 1.036723186076547, 1.036556299440207, 1.036013755642619, 1.035672659295609, 1.035853478369702, 1.036565395546309, 1.037391591092997)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $colJ[$i]
}

$colK = @(1.043287811949982, 1.043759798486559, 1.04406509757265, 1.044193418217483, 1.044214962192209, 1.044066812307455, 1.043447343546072, 1.042354992262529, 1.041626318711643, 1.041310707083756, 1.041193462483106, 1.041218612387583, 1.041301015853666, 1.041351785746192, 1.041647263006705, 1.041832584599774, 1.041940670665103, 1.041977523673429, 1.041812702256235, 1.041276750418258, 1.040939705192272, 1.041118385539897, 1.041821686259071, 1.042637474342081)
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Cells.Item($i + 2, 11).Value = $colK[$i]
}

$colL = @(1.035110029980152, 1.035626089392669, 1.03596027291138, 1.036100824204312, 1.036124426895476, 1.035962150728456, 1.035284380317929, 1.034092102778116, 1.033298704219001, 1.032955515699087, 1.032828095416812, 1.032855424979092, 1.032944981964437, 1.033000168336744, 1.033321488177851, 1.033523140527904, 1.033640795409096, 1.033680918488872, 1.033501501563041, 1.03291860812909, 1.03255244003156, 1.032746521860256, 1.033511279173473, 1.034400084920395)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colM = @(1.05196360234141, 1.052669151238456, 1.053126178818025, 1.053318428937666, 1.053350715317075, 1.053128747220874, 1.052201943392514, 1.05057262428407, 1.049489083248912, 1.049020552435653, 1.048846618497642, 1.04888392339887, 1.049006172967383, 1.049081508165199, 1.049520191925006, 1.049795541754989, 1.049956211217582, 1.050011005913585, 1.049765992840635, 1.04897017076275, 1.048470380995077, 1.048735273930346, 1.049779344530006, 1.050993378617811)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}

$colN = @(1.03952471153555, 1.040003837872065, 1.040313846328752, 1.040444168243487, 1.040466049499202, 1.040315587720291, 1.039686637563479, 1.038578262198921, 1.037839363084956, 1.037519431615907, 1.037400598604337, 1.037426088497239, 1.037509608757632, 1.037561068885547, 1.037860596343495, 1.038048487377046, 1.038158082566506, 1.038195451914587, 1.038028328280082, 1.03748501400807, 1.037143433265123, 1.037324509123042, 1.038037437303697, 1.038864806142828)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $colN[$i]
}
